$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be inserted at column B for each data row (row -> new leading value).
$newValues = @{
    2  = -0.5616080510579985
    3  = 0.1417647591280393
    4  = -0.4790798465348092
    5  = 0.1916007792754515
    6  = 1.573432754301089
    7  = 0.9422837133007778
    8  = 0.0678490295623069
    9  = -0.5264228954459207
    10 = 0.8949500190880419
    11 = 0.2303995154407018
    12 = 0.4008418571243615
    13 = 0.2679782848922332
    14 = -0.5417707991668423
    15 = 0.0506862842519193
    16 = -0.1624199859130616
}

# For each row, read existing values from columns B..K, shift them one
# column to the right (dropping anything that would overflow past column
# K), then insert the new value into column B.
for ($row = 2; $row -le 16; $row++) {
    $existing = @()
    for ($col = 2; $col -le 11; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value2
        if ($val -ne $null) {
            $existing += $val
        }
    }

    $combined = @($newValues[$row]) + $existing
    if ($combined.Count -gt 10) {
        $combined = $combined[0..9]
    }

    for ($i = 0; $i -lt 10; $i++) {
        $col = 2 + $i
        $cell = $ws.Cells.Item($row, $col)
        if ($i -lt $combined.Count) {
            $cell.Value = $combined[$i]
        } else {
            $cell.Value = $null
        }
    }
}
